$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new ticker entry in the next empty row (A5)
$ws.Range("A5").Value = "GRT-USD"
